# Update the BTC price history sheet: append the 2025 year-end price,
# matching the commit "update metrics and price data files".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Append the new data row (year 2025 / end-of-year BTC price).
$ws.Range("A18").Value = 2025
$ws.Range("B18").Value = 75269.850000000006

# Leave the selection where a user would land after typing the last
# value and pressing Enter/Tab past the edited row.
$ws.Range("D19").Select() | Out-Null
